$d = $word.ActiveDocument

function Get-TargetParagraphRange($d) {
    $rng = $d.Content
    $rng.Find.Execute("totalamountdue", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $para = $rng.Paragraphs(1)
    return $para.Range
}

$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:tabs><w:tab w:val="left" w:pos="5575"/></w:tabs><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="BC Sans" w:hAnsi="BC Sans" w:cs="BC Sans"/><w:b/><w:color w:val="0A3266"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr></w:p>'

# Insert first new empty paragraph right before the target paragraph
$pRange = Get-TargetParagraphRange $d
$insertPoint = $d.Range($pRange.Start, $pRange.Start)
$insertPoint.InsertXML($newParaXml)

# Insert second new empty paragraph right before the target paragraph (re-find since offsets shifted)
$pRange2 = Get-TargetParagraphRange $d
$insertPoint2 = $d.Range($pRange2.Start, $pRange2.Start)
$insertPoint2.InsertXML($newParaXml)

Write-Host "Paragraphs after insert: $($d.Paragraphs.Count)"

# Now split the run ".feedata.totalamountdue" into ".feedata." + new run "depositpaid"
$findRng = $d.Content
$findRng.Find.Execute("totalamountdue", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$s = $findRng.Start
$e = $findRng.End

$findRng.Text = ""
$ip = $d.Range($s, $s)
$ip.InsertAfter("depositpaid")

# force a run split by nudging formatting on just the new text, then restoring it,
# so the new text doesn't get silently re-merged into the preceding run
$newTextRng = $d.Range($s, $s + 11)
$newTextRng.Font.Size = 8
$newTextRng.Font.Size = 9

Write-Host "Done replacing totalamountdue -> depositpaid"
